$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add new header columns H1 (Hint) and I1 (Popup) ---
$ws.Range("H1").Value = "Hint"
$ws.Range("I1").Value = "Popup"

# Give the new header cells the same font/fill treatment as the rest of the
# header row (white text on the blue header fill), via the workbook's
# existing "MetaData_Columns" cell style (same font/fill as A1:G1, but
# without the border/number-format/wrap used on A1:G1).
$ws.Range("H1:I1").Style = "MetaData_Columns"

# --- Update E3 with the new multi-line Variables text, and enable wrap ---
$ws.Range("E3").Value = "item_api_category = TO_TEXT(SELECT a.CategoryChoice FROM Products p JOIN API a ON p.APIID == a.ID WHERE p.ID == [item]);`nitem_med_form = TO_TEXT(SELECT MedFormID FROM Products WHERE ID == [item]); `nSAVE(item_api_category); `nSAVE(item_med_form); `nGO(product);"
$ws.Range("E3").WrapText = $true

# Row 3 needs a custom height to fit the wrapped multi-line text
$ws.Rows.Item(3).RowHeight = 65.55

# --- Update selection to mirror the authored workbook state ---
$ws.Range("E10").Select()
